$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d7ad12ee08e9815c6c1c7f2c20b9557c833dc6e/e2e/f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ee1eef43992e86949afa05e4775ca18a8f15608/e2e/f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d7ad12ee08e9815c6c1c7f2c20b9557c833dc6e/e2e/f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.md."

# --- zh-cn sheet, row 7 (f128176a file) ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = "f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.md"
$wsZh.Range("J7").Value = "f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.c2fdf99cfe13beded97e5348e43174fc1e04054a.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-26 22:54:17"
$wsZh.Range("P7").Value = $errorDetail

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", "f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.md")
$wsZh.Range("I7").Font.Underline = $true
$wsZh.Range("I7").Font.Color = 15570276

# --- de-de sheet, row 7 (f128176a file) ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = "f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.md"
$wsDe.Range("J7").Value = "f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.c2fdf99cfe13beded97e5348e43174fc1e04054a.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-26 22:54:24"
$wsDe.Range("P7").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", "f128176a-8a0c-4f9a-bbf9-b762ec22cfbf.md")
$wsDe.Range("I7").Font.Underline = $true
$wsDe.Range("I7").Font.Color = 15570276
